{"js": "// Fix templates 002/005: the liquidation paragraph referred to the\n// shareholder's address fields with a stray \"sh.\" prefix\n// (e.g. \"{{ sh.street_number }}\") instead of the full\n// \"{{ shareholders[0]. ... }}\" accessor used elsewhere in the template.\n// Replace every \"sh.\" occurrence with \"shareholders[0].\" so the four\n// address placeholders (street_number, street_name, zip_code, city)\n// resolve correctly.\n\nconst body = context.document.body;\n\n// Find every occurrence of the stray \"sh.\" prefix in the document body.\nconst searchResults = body.search(\"sh.\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\n// Walk backwards so earlier ranges stay valid while later ones are edited.\nfor (let i = searchResults.items.length - 1; i >= 0; i--) {\n  searchResults.items[i].insertText(\"shareholders[0].\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Fix templates 002/005: the liquidation paragraph referred to the\n# shareholder's address fields with a stray \"sh.\" prefix\n# (e.g. \"{{ sh.street_number }}\") instead of the full\n# \"{{ shareholders[0]. ... }}\" accessor used elsewhere in the template.\n# Replace every \"sh.\" occurrence with \"shareholders[0].\" so the four\n# address placeholders (street_number, street_name, zip_code, city)\n# resolve correctly.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"sh.\"\n$find.Replacement.Text = \"shareholders[0].\"\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n#   ReplaceWith, Replace)\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\n    $find.Text,\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find.Replacement.Text,\n    2\n) | Out-Null\n"}
